# Add a second quiz/worksheet-score column pair (Q2/W2) next to the existing
# Q1/W1 columns, populate it for every existing student, insert a new student
# row (Sturdifen, Jasmine) before Thomas/Troetti, and tidy up the page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: I1 = "Q2", J1 = "W2" ---------------------------------
$ws.Range("I1").Value = "Q2"
$ws.Range("J1").Value = "W2"

# --- New Q2/W2 scores for the existing student rows (rows 3-23) ------------
# Row: [Q2 (I), W2 (J)]
$scores = @{
    3  = @(4, 50)
    4  = @(4, 50)
    5  = @(4, 50)
    6  = @(4, 47)
    7  = @(4, 50)
    8  = @(4, 37)
    9  = @(4, 48)
    10 = @(4, 47)
    11 = @(4, 50)
    12 = @(4, 50)
    13 = @(2, 50)
    14 = @(4, 50)
    15 = @(4, 47)
    16 = @(4, 47)
    17 = @(4, 50)
    18 = @(4, 50)
    19 = @(4, 50)
    20 = @(4, 47)
    21 = @(2, 50)
    22 = @(2, 50)
    23 = @(4, 47)
}

foreach ($r in $scores.Keys) {
    $pair = $scores[$r]
    $ws.Range("I$r").Value = $pair[0]
    $ws.Range("J$r").Value = $pair[1]
}

# --- Insert a new student row at row 24 (Thomas/Troetti shift down) --------
$ws.Rows.Item(24).Insert() | Out-Null

$ws.Range("A24").Value = "Sturdifen"
$ws.Range("B24").Value = "Jasmine"
$ws.Range("G24").Value = "?"
$ws.Range("H24").Value = "?"
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 50

# --- Q2/W2 scores for the rows that shifted down (now 25 and 26) -----------
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 50
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = 47

# --- Cosmetic bits: page orientation + last active selection ---------------
$ws.PageSetup.Orientation = 1
$ws.Range("J7").Select() | Out-Null
